# Generate Report for Handoff
# Updates the "b.md" rows across the Overview / zh-cn / de-de sheets to
# reflect that a new handoff package has been generated for b.md.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3160f06d2de5ef52c8855b025e33a6918873c8fc/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2afe580469ec78173171a53f4067387aeb4bba15/e2e/b.md."

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-07 13:15:38"

# --- zh-cn sheet --------------------------------------------------------
# (values are prefixed with a leading apostrophe so Excel stores them as
# plain text, matching the original workbook where "True"/"False" are
# shared text strings rather than boolean cells)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-07 13:15:24"
$wsZhCn.Range("P3").Value = $errorDetail
# Column G already has the target width (40 "OOXML" units); reuse its
# ColumnWidth so the conversion Excel applies lines up exactly.
$wsZhCn.Columns.Item(16).ColumnWidth = $wsZhCn.Columns.Item(7).ColumnWidth

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-07 13:15:38"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $wsDeDe.Columns.Item(7).ColumnWidth
